$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.500.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.442.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.441.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.887.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000168"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.355.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.451.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "318.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0968"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.565.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "526.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.51%  "

$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.74%  "

$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "138.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.50%  "

$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("E45").Value = "  -2.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0524"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.586"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0930"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
